$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New footer mapping sample data (mirrors "offset_footer" field in the
# account_statement_import_txt_xlsx mapping test fixture).
$ws.Range("C6").Value = "Any footer data"
$ws.Range("C6").WrapText = $true
$ws.Range("E6").WrapText = $true

$ws.Range("E7").Value = "Any footer data"
$ws.Range("E7").WrapText = $true

$ws.Range("F8").Value = "Any footer data"
$ws.Range("F8").WrapText = $true

# Rows with wrapped text grow taller to show the wrapped content.
$ws.Rows.Item(6).RowHeight = 23.85
$ws.Rows.Item(8).RowHeight = 23.85

# Match the selection left behind by the author when they made the edit.
$selected = $ws.Range("E7").Select()
